$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (the "ECs" target-cluster row); rows 3 and 4 shift up to become 2 and 3
$ws.Rows("2").Delete()

# Update recalculated values on the new row 2 (previously row 3, FAPs target)
$ws.Range("G2").Value = 0.03588566666666667
$ws.Range("O2").Value = 0.6232066589444157
$ws.Range("P2").Value = 0.6232066589444157
$ws.Range("S2").Value = 0.6232066589444157
$ws.Range("T2").Value = 0.6232066589444157

# Update recalculated values on the new row 3 (previously row 4, MuSCs target)
$ws.Range("G3").Value = 0.03588566666666667
$ws.Range("O3").Value = 0.3767933410555842
$ws.Range("P3").Value = 0.3767933410555843
$ws.Range("S3").Value = 0.3767933410555842
$ws.Range("T3").Value = 0.3767933410555843
